# Auto-applied numeric updates to the Leve profit tables (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1500374.8
$ws.Range("J17").Value = 1552108.4
$ws.Range("L17").Value = 4656325.199999999
$ws.Range("N17").Value = -4656661.199999999
$ws.Range("H69").Value = 5622.5
$ws.Range("I69").Value = 5245
$ws.Range("K69").Value = 15735
$ws.Range("M69").Value = -14861
$ws.Range("H72").Value = 5622.5
$ws.Range("I72").Value = 5245
$ws.Range("K72").Value = 47205
$ws.Range("M72").Value = -42837
$ws.Range("H113").Value = 5074.8887
$ws.Range("I113").Value = 5796.25
$ws.Range("K113").Value = 5796.25
$ws.Range("M113").Value = -2542.25
$ws.Range("H132").Value = 4878.304
$ws.Range("I132").Value = 5273.2104
$ws.Range("J132").Value = 3002.5
$ws.Range("K132").Value = 15819.6312
$ws.Range("L132").Value = 9007.5
$ws.Range("M132").Value = -13289.6312
$ws.Range("N132").Value = -14067.5
$ws.Range("H135").Value = 1167.8
$ws.Range("I135").Value = 1130.8889
$ws.Range("J135").Value = 1500
$ws.Range("K135").Value = 10178.0001
$ws.Range("L135").Value = 13500
$ws.Range("M135").Value = -7643.000099999999
$ws.Range("N135").Value = -18570
$ws.Range("H137").Value = 17705
$ws.Range("I137").Value = 26489.691
$ws.Range("K137").Value = 79469.073
$ws.Range("M137").Value = -76919.073

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3460.8813
$ws.Range("J32").Value = 10949.5
$ws.Range("L32").Value = 10949.5
$ws.Range("N32").Value = -11523.5
$ws.Range("H132").Value = 70183.03
$ws.Range("I132").Value = 3699.6667
$ws.Range("K132").Value = 11099.0001
$ws.Range("M132").Value = -8569.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 39373.375
$ws.Range("I75").Value = 10747
$ws.Range("J75").Value = 67999.75
$ws.Range("K75").Value = 10747
$ws.Range("L75").Value = 67999.75
$ws.Range("M75").Value = -9811
$ws.Range("N75").Value = -69871.75
$ws.Range("H78").Value = 39373.375
$ws.Range("I78").Value = 10747
$ws.Range("J78").Value = 67999.75
$ws.Range("K78").Value = 32241
$ws.Range("L78").Value = 203999.25
$ws.Range("M78").Value = -27561
$ws.Range("N78").Value = -213359.25
$ws.Range("H134").Value = 1440.8572
$ws.Range("I134").Value = 1440.8572
$ws.Range("K134").Value = 4322.571599999999
$ws.Range("M134").Value = -1787.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 346666660
$ws.Range("I4").Value = 20000000
$ws.Range("J4").Value = 1000000000
$ws.Range("K4").Value = 20000000
$ws.Range("L4").Value = 1000000000
$ws.Range("M4").Value = -19999888
$ws.Range("N4").Value = -1000000224
$ws.Range("H10").Value = 500546.75
$ws.Range("I10").Value = 667062.3
$ws.Range("K10").Value = 667062.3
$ws.Range("M10").Value = -666923.3
$ws.Range("H16").Value = 2306.7144
$ws.Range("I16").Value = 2262.182
$ws.Range("J16").Value = 2470
$ws.Range("K16").Value = 2262.182
$ws.Range("L16").Value = 2470
$ws.Range("M16").Value = -1975.182
$ws.Range("N16").Value = -3044
$ws.Range("H31").Value = 20838200
$ws.Range("I31").Value = 2938
$ws.Range("J31").Value = 35720530
$ws.Range("K31").Value = 2938
$ws.Range("L31").Value = 35720530
$ws.Range("M31").Value = -2643
$ws.Range("N31").Value = -35721120
$ws.Range("H33").Value = 2009.6
$ws.Range("I33").Value = 2009.6
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 2009.6
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -1630.6
$ws.Range("N33").ClearContents()
$ws.Range("H34").Value = 20838200
$ws.Range("I34").Value = 2938
$ws.Range("J34").Value = 35720530
$ws.Range("K34").Value = 2938
$ws.Range("L34").Value = 35720530
$ws.Range("M34").Value = -2736
$ws.Range("N34").Value = -35720934
$ws.Range("H45").Value = 9999
$ws.Range("I45").Value = 9999
$ws.Range("K45").Value = 9999
$ws.Range("M45").Value = -9406
$ws.Range("H58").Value = 19645594
$ws.Range("I58").Value = 3155.6365
$ws.Range("K58").Value = 3155.6365
$ws.Range("M58").Value = -2952.6365
$ws.Range("H94").Value = 2634.875
$ws.Range("J94").Value = 2768.077
$ws.Range("L94").Value = 2768.077
$ws.Range("N94").Value = -3670.077
$ws.Range("H107").Value = 671.1429000000001
$ws.Range("J107").Value = 715.6667
$ws.Range("L107").Value = 715.6667
$ws.Range("N107").Value = -4555.6667
$ws.Range("H113").Value = 2306.7144
$ws.Range("I113").Value = 2262.182
$ws.Range("J113").Value = 2470
$ws.Range("K113").Value = 2262.182
$ws.Range("L113").Value = 2470
$ws.Range("M113").Value = -92.18199999999979
$ws.Range("N113").Value = -6810
$ws.Range("H134").Value = 4089.3125
$ws.Range("I134").Value = 3369.3333
$ws.Range("J134").Value = 6249.25
$ws.Range("K134").Value = 10107.9999
$ws.Range("L134").Value = 18747.75
$ws.Range("M134").Value = -7572.999899999999
$ws.Range("N134").Value = -23817.75
$ws.Range("H136").Value = 19645594
$ws.Range("I136").Value = 3155.6365
$ws.Range("K136").Value = 9466.9095
$ws.Range("M136").Value = -6916.9095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 1297.7778
$ws.Range("I136").Value = 1297.7778
$ws.Range("K136").Value = 3893.3334
$ws.Range("M136").Value = 1206.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4170.037
$ws.Range("I102").Value = 1824.2
$ws.Range("J102").Value = 7102.3335
$ws.Range("K102").Value = 1824.2
$ws.Range("L102").Value = 7102.3335
$ws.Range("M102").Value = -202.2
$ws.Range("N102").Value = -10346.3335
$ws.Range("H132").Value = 2314.9412
$ws.Range("I132").Value = 1709.625
$ws.Range("K132").Value = 5128.875
$ws.Range("M132").Value = -2598.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4170774.2
$ws.Range("I122").Value = 4042.2273
$ws.Range("K122").Value = 12126.6819
$ws.Range("M122").Value = -9676.6819
$ws.Range("H132").Value = 2002.3462
$ws.Range("I132").Value = 1485.3529
$ws.Range("K132").Value = 4456.0587
$ws.Range("M132").Value = -1926.0587

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 25499
$ws.Range("J45").Value = 25499
$ws.Range("L45").Value = 25499
$ws.Range("N45").Value = -26481
$ws.Range("H62").Value = 8766
$ws.Range("I62").Value = 5997.5
$ws.Range("J62").Value = 9557
$ws.Range("K62").Value = 5997.5
$ws.Range("L62").Value = 9557
$ws.Range("M62").Value = -5373.5
$ws.Range("N62").Value = -10805
$ws.Range("H65").Value = 8766
$ws.Range("I65").Value = 5997.5
$ws.Range("J65").Value = 9557
$ws.Range("K65").Value = 29987.5
$ws.Range("L65").Value = 47785
$ws.Range("M65").Value = -26867.5
$ws.Range("N65").Value = -54025
$ws.Range("H107").Value = 1329.4348
$ws.Range("I107").Value = 1054.4445
$ws.Range("K107").Value = 3163.3335
$ws.Range("M107").Value = -1243.3335
$ws.Range("H132").Value = 2906
$ws.Range("I132").Value = 2742.3823
$ws.Range("K132").Value = 8227.1469
$ws.Range("M132").Value = -5697.1469
